$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the indicator values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 147498.42000000001
$metrics.Range("B3").Value  = 129247.56
$metrics.Range("B4").Value  = 45997.86
$metrics.Range("B5").Value  = 6120
$metrics.Range("B6").Value  = 4943744.1700000009
$metrics.Range("B7").Value  = 4171324.2400000007
$metrics.Range("B8").Value  = 1452957.69
$metrics.Range("B9").Value  = 192327
$metrics.Range("B10").Value = 33409125.160000008
$metrics.Range("B11").Value = 31446599.399999999
$metrics.Range("B12").Value = 11734679.730000002
$metrics.Range("B13").Value = 1289957

# Move the Metrics sheet's selection to D16 (matches recorded selection change)
$metrics.Range("D16").Select()

# --- "today" sheet: move its selection to E6. Its B11:B22/E11:E22/F11:F22 ---
# --- cells are formulas referencing Metrics! so they recalc automatically. ---
$today = $wb.Worksheets.Item("today")
$today.Range("E6").Select()
